$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Stage current D:G contents (columns 4-7) into helper columns Z:AC (26-29)
# so the text/number cell-type (shared string, no quote-prefix) is preserved
# via Copy/PasteSpecial instead of a literal Value2 assignment (which would
# coerce numeric-looking strings like "110" into real numbers).
$srcRange = $ws.Range($ws.Cells(1, 4), $ws.Cells($lastRow, 7))
$srcRange.Copy()
$stageRange = $ws.Range($ws.Cells(1, 26), $ws.Cells($lastRow, 29))
$stageRange.PasteSpecial(-4163)

for ($i = 1; $i -le $lastRow; $i++) {
    # new D = old G (staged col AC = 29)
    $ws.Cells($i, 29).Copy()
    $ws.Cells($i, 4).PasteSpecial(-4163)

    # new E = old F (staged col AB = 28)
    $ws.Cells($i, 28).Copy()
    $ws.Cells($i, 5).PasteSpecial(-4163)

    # new F = old D (staged col Z = 26)
    $ws.Cells($i, 26).Copy()
    $ws.Cells($i, 6).PasteSpecial(-4163)

    # new G = old E (staged col AA = 27)
    $ws.Cells($i, 27).Copy()
    $ws.Cells($i, 7).PasteSpecial(-4163)
}

# Clean up the staging columns
$stageRange.ClearContents()
$excel.CutCopyMode = 0
